$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '302.51'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2.17%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '44.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '6.74%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.097'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.18%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07706'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.19%'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.617'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.97%'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.047'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '12.87%'
$ws.Range('B8').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C8').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.1273'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '7.06%'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1867'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '1.92%'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09202'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '3.02%'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.04159'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.48%'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1047'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.57%'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001276'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.49%'
$ws.Range('B14').Value = 'TigerCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.005768'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-2.11%'
$ws.Range('B15').Value = 'UpBots'
$ws.Range('C15').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.007489'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1,911.15%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.346'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.02%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.416'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.39%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-3.19%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3343'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.93%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.104'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '2.39%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1398'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.75%'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '7.16%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04191'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '3.69%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001282'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.37%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004420'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '14.24%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001350'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '9.72%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02500'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '4.27%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05295'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '1.88%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.005934'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-9.81%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007733'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.73%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1350'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '2.16%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007346'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-0.30%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007550'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '2.45%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3013'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-6.42%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006673'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '7.39%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.06%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.04336'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-4.78%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.01%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002099'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.06%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0001999'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.06%'
